$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new item ("قصافات كبار") was added to the nواقص (shortages) list, in
# its correct alphabetical spot just before "قطن 50جم". On this report the
# "م" counter column (A) is positional/static, so the cleanest way to
# reproduce the generated file is to:
#   1) open up a fresh row right before the current totals row (106),
#      which leaves rows 103-105 and their A/B numbering untouched and
#      pushes the totals/footer rows down to 107/108;
#   2) push the existing item rows 103-105 down one VISUAL item slot by
#      re-writing their item columns (C/H/N/P/Q) with the next row's old
#      data, cascading from the bottom up;
#   3) drop the new item into row 103;
#   4) fill the freshly inserted row 106 with what used to be the last
#      item (previously on row 105);
#   5) bump the grand total and the generated timestamp.

# 1) Make room for one more item row just above the totals row.
$ws.Rows("106:106").Insert()

# Give the new row the same look as a normal item row (format only - the
# values are overwritten below).
$ws.Range("A105:Q105").Copy()
$ws.Range("A106:Q106").PasteSpecial(-4122)
$ws.Rows("106:106").RowHeight = 25.5
$ws.Range("A106:B106").Merge()
$ws.Range("C106:G106").Merge()
$ws.Range("H106:K106").Merge()
$ws.Range("L106:M106").Merge()
$ws.Range("N106:O106").Merge()

# 2) Cascade the item data down one slot, starting from the bottom so we
#    never overwrite a row before we've read its old value. (`.Value2` is
#    used everywhere instead of `.Value` for reading, since this host's
#    `.Value` getter doesn't resolve properly.)
$ws.Range("C106").Value2 = $ws.Range("C105").Value2
$ws.Range("H106").Value2 = $ws.Range("H105").Value2
$ws.Range("N106").Value2 = $ws.Range("N105").Value2
$ws.Range("P106").Value2 = $ws.Range("P105").Value2
$ws.Range("Q106").Value2 = $ws.Range("Q105").Value2
$ws.Range("A106").Value2 = 100

$ws.Range("C105").Value2 = $ws.Range("C104").Value2
$ws.Range("H105").Value2 = $ws.Range("H104").Value2
$ws.Range("N105").Value2 = $ws.Range("N104").Value2
$ws.Range("P105").Value2 = $ws.Range("P104").Value2

$ws.Range("C104").Value2 = $ws.Range("C103").Value2
$ws.Range("H104").Value2 = $ws.Range("H103").Value2
$ws.Range("N104").Value2 = $ws.Range("N103").Value2
$ws.Range("P104").Value2 = $ws.Range("P103").Value2

# 3) New item goes on row 103.
$ws.Range("C103").Value2 = "قصافات كبار"
$ws.Range("H103").Value2 = "9:0"
$ws.Range("N103").Value2 = "20.00"
$ws.Range("P103").Value2 = "20.0000"

# 5) Grand total (now row 107) grows by the new item's amount, and the
#    generated-on timestamp (now row 108) is refreshed.
$ws.Range("P107").Value2 = 5933.795
$ws.Range("A108").Value2 = "Friday, 18 July, 2025 10:55 PM"
